# TradesmenDirectory.pptx edit
# Commit: "updated the ppt, I added a slide on the phpBB software thet
#          we're going to use for the jobs board."
#
# Two content changes are applied:
#  1. Slide 3 ("Features"): the "Email " / "verification" runs are
#     merged back into a single "Email verification" run.
#  2. Slide 6 ("Jobs Board"): the content placeholder gains three new
#     paragraphs describing the phpBB software used for the jobs board,
#     and the existing sentence is re-split across a few runs.

$p = $ppt.ActivePresentation

# --- Slide 3: "Features" -------------------------------------------------
$s3 = $p.Slides.Item(3)
$content3 = $s3.Shapes.Item(1)
$tr3 = $content3.TextFrame.TextRange

# Locate the "Email verification" paragraph (currently split into the
# runs "Email " + "verification") and rewrite it as a single run,
# leaving every other paragraph on the slide untouched. It is the 4th
# paragraph in the placeholder: "Consumer feedback" (18 chars) + break,
# "Current jobs being completed" (29 chars) + break, "Verified
# tradesmen / waiting period" (36 chars) + break -> absolute start 84,
# length 19 ("Email verification").
$tr3.Characters(84, 19).Text = "Email verification"

# --- Slide 6: "Jobs Board" ------------------------------------------------
$s6 = $p.Slides.Item(6)
$content6 = $s6.Shapes.Item(1)
$tr6 = $content6.TextFrame.TextRange

$tr6.Text = "Verified consumers will have the ability to advertise jobs on the job board, where tradesmen can reply back to them.`rThe jobs board is created using a software called phpBB.`rphpBB is a widely used open source bulletin board system in the world.`rIt supports private messaging. Which we can use to allow consumers and tradesmen to message each other on the site."

# Re-split paragraph 1 into its four runs.
$tr6.Characters(1, 77).Text = "Verified consumers will have the ability to advertise jobs on the job board, "
$tr6.Characters(78, 6).Text = "where "
$tr6.Characters(84, 32).Text = "tradesmen can reply back to them"
$tr6.Characters(116, 1).Text = "."

# Re-split paragraph 3 (the "phpBB is a widely used ..." sentence) into
# its five runs. Paragraph 3 starts at absolute character 175 (116 chars
# of paragraph 1 + 1 paragraph break + 56 chars of paragraph 2 + 1 break).
$tr6.Characters(175, 5).Text = "phpBB"
$tr6.Characters(180, 4).Text = " is "
$tr6.Characters(184, 2).Text = "a "
$tr6.Characters(186, 53).Text = "widely used open source bulletin board system in the "
$tr6.Characters(239, 6).Text = "world."
